$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.532399296760559
$ws.Range("B1").Value = 2.161391258239746
$ws.Range("C1").Value = 2.632901430130005
$ws.Range("D1").Value = 4.127182483673096
$ws.Range("E1").Value = 0.6179674863815308
